# Apply weekly price-update shuffle for "Vega Modelo de Temuco - Pera asiática"
# Rows 2-11 of the sheet get their D,L,M,N,O,P,Q,R,S,T values re-assigned from
# a different row of the original data (a rotation of two cycles); row 9 is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values (before any mutation) for the relevant columns.
$cols = @("D","L","M","N","O","P","Q","R","S","T")
$orig = @{}
for ($r = 2; $r -le 11; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $row
}

# Mapping: target row -> source row (data copied from source row's original values)
$mapping = @{
    2  = 10
    3  = 7
    4  = 5
    5  = 6
    6  = 2
    7  = 11
    8  = 3
    9  = 9
    10 = 4
    11 = 8
}

foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    $srcData = $orig[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value2 = $srcData[$c]
    }
}
